# Apply the "include new vars in 2_recode" edit:
# Rows 17-20 and 24-26 in Sheet1 get updated variables_fuente (C) and situacion (E)
# values, shifting the perper_p_delito_pronostico_* numbering / labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = "perper_p_delito_pronostico_2"
$ws.Range("E17").Value = "Robo o hurto de su vehículo o portonazo"

$ws.Range("C18").Value = "perper_p_delito_pronostico_3"
$ws.Range("E18").Value = "Robo o hurto de algún objeto dejado dentro del vehículo o parte de él"

$ws.Range("C19").Value = "perper_p_delito_pronostico_4"
$ws.Range("E19").Value = "Vandalismo o daño a su vivienda o vehículo"

$ws.Range("C20").Value = "perper_p_delito_pronostico_6"
$ws.Range("E20").Value = "Hurto"

$ws.Range("C24").Value = "perper_p_delito_pronostico_5"
$ws.Range("E24").Value = "Robo o asalto, como robo con violencia, cogoteo, robo por sorpresa o lanzazo"

$ws.Range("C25").Value = "perper_p_delito_pronostico_7"
$ws.Range("E25").Value = "Agresiones físicas o lesiones"

$ws.Range("C26").Value = "perper_p_delito_pronostico_8"
$ws.Range("E26").Value = "Amenazas o extorsión"
